$wb = $excel.ActiveWorkbook

# Remove the extra "Now" worksheet (sheetId=4) that was added for staging/testing.
$nowSheet = $wb.Worksheets.Item("Now")
$nowSheet.Delete()

# Clean up the product list on Plan1: insert a new row for a locally-coded
# item (internal code 665544) ahead of the remaining GTIN-coded rows, adding
# a short clarifying comment via the inserted row.
$ws = $wb.Worksheets.Item("Plan1")
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 665544

# Renumber the trailing id column so the sequence stays contiguous after the insert.
for ($r = 7; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update the active selection to reflect where editing left off.
$ws.Range("C17").Select()
